$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 2 ("Поставлена задача: ...RELY.") - merge adjacent runs that
#    share identical formatting into single <w:r> elements (two places):
#      a) "Поставлена задача: разработать модель " + "бизнес-процесса" + " "
#      b) "управления рассылкой" + " "                     (bold run)
#    Leave the "для программного комплекса " run untouched in between.
# ---------------------------------------------------------------------------
$para2 = $d.Paragraphs(2).Range
$frag2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="007F3255" w14:textId="2A74DB12" w:rsidR="006B326B" w:rsidRPr="000E2E0B" w:rsidRDefault="000E2E0B" w:rsidP="000E2E0B" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Поставлена задача: разработать модель бизнес-процесса </w:t></w:r><w:r w:rsidRPr="00CC1F90"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">для программного комплекса </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">управления рассылкой </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>RELY</w:t></w:r><w:r w:rsidRPr="000E2E0B"><w:rPr><w:b/></w:rPr><w:t>.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$para2.InsertXML($frag2)

# ---------------------------------------------------------------------------
# 2) Relocate the "_GoBack" bookmark: remove it from around the diagram
#    picture and place a fresh (collapsed) copy right at the end of the
#    "Серверная часть принимает запрос ... в очередь." list item.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$para7 = $d.Paragraphs(7).Range
$frag7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="4BFDC4D2" w14:textId="482A40AA" w:rsidR="000E2E0B" w:rsidRDefault="000E2E0B" w:rsidP="000E2E0B" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Серверная часть принимает запрос от клиентской части, формирует письмо, заносит его в хранилище данных и отправляет письмо в очередь.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$para7.InsertXML($frag7)

# ---------------------------------------------------------------------------
# 3) Mark the run that holds the BPMN diagram picture as "no proofing" so it
#    gains <w:rPr><w:noProof/></w:rPr>.
# ---------------------------------------------------------------------------
$pictureParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ppr = $d.Paragraphs($i)
    if ($ppr.Range.InlineShapes.Count -gt 0) {
        $pictureParaIndex = $i
    }
}
$picturePara = $d.Paragraphs($pictureParaIndex).Range
$picturePara.NoProofing = 1
